$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new row at 17 (pushes the signature block rows 21/22 down to 22/23)
$ws.Rows("17:17").Insert()

# Copy formatting of row 16 (B:J) into the new row 17 (same style as the existing data row)
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the new row's data - duplicate worker record for the new period "2508"
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1002324705"
$ws.Range("D17").Value = "RONAL JOSE NAVARRO SARAVIA"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Update the totals to reflect the second period
$ws.Range("E11").Value = 113880
$ws.Range("F13").Value = 2
